$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: add Name_Regression value (column A) ---
$ws.Range("A20").Value2 = "daily_emissions"

# --- Row 21: daily_emissions_noattitude ---
$ws.Range("A21").Value2 = "daily_emissions_noattitude"
$ws.Range("B21").Value2 = "emissions_wege_wout_work"
$ws.Range("C2").Copy($ws.Range("C21"))                 # "FALSE" as text (shared string), not boolean
$ws.Range("E21").Value2 = "Regression for daily emissions, without attitudinal variables"
$ws.Range("F21").Value2 = 1
$ws.Range("G21").Value2 = "NoNA"
$ws.Range("H21").Value2 = 1
$ws.Range("I21").Value2 = 1
$ws.Range("J21").Value2 = 0

# --- Row 22: daily_emissions_spatial (note: E22 string registered before A22 string) ---
$ws.Range("E22").Value2 = "Regression for daily emissions, with additional controls for spatial characteristics"
$ws.Range("A22").Value2 = "daily_emissions_spatial"
$ws.Range("B22").Value2 = "emissions_wege_wout_work"
$ws.Range("C2").Copy($ws.Range("C22"))                 # "FALSE" as text
$ws.Range("F22").Value2 = 1
$ws.Range("G22").Value2 = "NoNA"
$ws.Range("H22").Value2 = 1
$ws.Range("I22").Value2 = 0
$ws.Range("J22").Value2 = 1

# --- Row 23: daily_emissions_withbusiness ---
$ws.Range("A23").Value2 = "daily_emissions_withbusiness"
$ws.Range("B23").Value2 = "emissions_wege"
$ws.Range("C2").Copy($ws.Range("C23"))                 # "FALSE" as text
$ws.Range("E23").Value2 = "Regression for daily emissions, including business travels. "
$ws.Range("F23").Value2 = 1
$ws.Range("G23").Value2 = "NoNA"
$ws.Range("H23").Value2 = 1
$ws.Range("I23").Value2 = 0
$ws.Range("J23").Value2 = 0

# --- Row 24: daily_emissions_control_frequency ---
$ws.Range("A24").Value2 = "daily_emissions_control_frequency"
$ws.Range("B24").Value2 = "emissions_wege_wout_work"
$ws.Range("C3").Copy($ws.Range("C24"))                 # "TRUE" as text
$ws.Range("D24").Value2 = "P_NUTZ_RAD,P_NUTZ_OPNV"
$ws.Range("E24").Value2 = "Regression for daily emissions, when adding a control variables for the use of low-carbon transportation modes"
$ws.Range("F24").Value2 = 1
$ws.Range("G24").Value2 = "NoNA"
$ws.Range("H24").Value2 = 1
$ws.Range("I24").Value2 = 0
$ws.Range("J24").Value2 = 0

$excel.CutCopyMode = $false

# --- Update the view: scroll position + selection to match the edited area ---
$ws.Range("D23").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 3 | Out-Null
